{"js": "// Replace the date heading and each \"dividend\u00f7divisor=quotient, remainder\"\n// cell value in the worksheet table with the new values from the commit.\nconst replacements = [\n  [\"2024-07-22 Monday\", \"2024-07-23 Tuesday\"],\n  [\"674\u00f75=134, 4\", \"470\u00f75=94, 0\"],\n  [\"413\u00f75=82, 3\", \"662\u00f75=132, 2\"],\n  [\"667\u00f77=95, 2\", \"654\u00f76=109, 0\"],\n  [\"743\u00f78=92, 7\", \"121\u00f79=13, 4\"],\n  [\"798\u00f76=133, 0\", \"351\u00f77=50, 1\"],\n  [\"209\u00f72=104, 1\", \"448\u00f76=74, 4\"],\n  [\"349\u00f79=38, 7\", \"523\u00f73=174, 1\"],\n  [\"935\u00f76=155, 5\", \"460\u00f77=65, 5\"],\n  [\"155\u00f76=25, 5\", \"430\u00f78=53, 6\"],\n  [\"658\u00f73=219, 1\", \"750\u00f79=83, 3\"],\n  [\"103\u00f73=34, 1\", \"911\u00f72=455, 1\"],\n  [\"342\u00f78=42, 6\", \"532\u00f73=177, 1\"],\n  [\"864\u00f75=172, 4\", \"578\u00f76=96, 2\"],\n  [\"379\u00f78=47, 3\", \"842\u00f76=140, 2\"],\n  [\"284\u00f73=94, 2\", \"980\u00f76=163, 2\"],\n  [\"787\u00f78=98, 3\", \"917\u00f76=152, 5\"],\n  [\"983\u00f76=163, 5\", \"479\u00f72=239, 1\"],\n  [\"806\u00f78=100, 6\", \"708\u00f79=78, 6\"],\n  [\"870\u00f78=108, 6\", \"375\u00f75=75, 0\"],\n  [\"701\u00f75=140, 1\", \"894\u00f78=111, 6\"],\n  [\"492\u00f76=82, 0\", \"183\u00f74=45, 3\"],\n  [\"188\u00f79=20, 8\", \"878\u00f73=292, 2\"],\n  [\"198\u00f76=33, 0\", \"467\u00f73=155, 2\"],\n  [\"583\u00f74=145, 3\", \"573\u00f79=63, 6\"],\n  [\"438\u00f74=109, 2\", \"600\u00f76=100, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date heading and each \"dividend\u00f7divisor=quotient, remainder\"\n# cell value in the worksheet table with the new values from the commit.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-07-22 Monday\", \"2024-07-23 Tuesday\"),\n    @(\"674\u00f75=134, 4\", \"470\u00f75=94, 0\"),\n    @(\"413\u00f75=82, 3\", \"662\u00f75=132, 2\"),\n    @(\"667\u00f77=95, 2\", \"654\u00f76=109, 0\"),\n    @(\"743\u00f78=92, 7\", \"121\u00f79=13, 4\"),\n    @(\"798\u00f76=133, 0\", \"351\u00f77=50, 1\"),\n    @(\"209\u00f72=104, 1\", \"448\u00f76=74, 4\"),\n    @(\"349\u00f79=38, 7\", \"523\u00f73=174, 1\"),\n    @(\"935\u00f76=155, 5\", \"460\u00f77=65, 5\"),\n    @(\"155\u00f76=25, 5\", \"430\u00f78=53, 6\"),\n    @(\"658\u00f73=219, 1\", \"750\u00f79=83, 3\"),\n    @(\"103\u00f73=34, 1\", \"911\u00f72=455, 1\"),\n    @(\"342\u00f78=42, 6\", \"532\u00f73=177, 1\"),\n    @(\"864\u00f75=172, 4\", \"578\u00f76=96, 2\"),\n    @(\"379\u00f78=47, 3\", \"842\u00f76=140, 2\"),\n    @(\"284\u00f73=94, 2\", \"980\u00f76=163, 2\"),\n    @(\"787\u00f78=98, 3\", \"917\u00f76=152, 5\"),\n    @(\"983\u00f76=163, 5\", \"479\u00f72=239, 1\"),\n    @(\"806\u00f78=100, 6\", \"708\u00f79=78, 6\"),\n    @(\"870\u00f78=108, 6\", \"375\u00f75=75, 0\"),\n    @(\"701\u00f75=140, 1\", \"894\u00f78=111, 6\"),\n    @(\"492\u00f76=82, 0\", \"183\u00f74=45, 3\"),\n    @(\"188\u00f79=20, 8\", \"878\u00f73=292, 2\"),\n    @(\"198\u00f76=33, 0\", \"467\u00f73=155, 2\"),\n    @(\"583\u00f74=145, 3\", \"573\u00f79=63, 6\"),\n    @(\"438\u00f74=109, 2\", \"600\u00f76=100, 0\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
